# config test data changes
# Update the email addresses (hyperlink cell text) for the three
# registered test users, and move the selection to C8 afterwards
# (matching the state Excel saved the workbook in after this edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "107automanickam@gmail.com"
$ws.Range("C3").Value = "107automonika@gmail.com"
$ws.Range("C4").Value = "107autoharitha@gmail.com"

$ws.Range("C8").Select()
